# Apply the two changes captured by the commit:
#   1. Slide 6's table switches to a different built-in PowerPoint table style.
#   2. The deck's theme colour palette moves from the "Integral" scheme to the
#      default "Office Theme" scheme.

function HexToRgbLong($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1) Table style change on slide 6's table ---
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{8C4F2756-433A-4418-AF57-2B60A546F68C}")

# --- 2) Theme colours: Integral -> Office Theme ---
# Order matches the standard theme colour slots:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToRgbLong($officeThemeColors[$i - 1])
}
